$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 200, shifting existing rows 200:208 down to 201:209.
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with the latest weekly price entry.
$ws.Cells.Item(200, 1).Value = 2
$ws.Cells.Item(200, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(200, 3).Value = "Coquimbo"
$ws.Cells.Item(200, 4).Value = "2023-04-05"
$ws.Cells.Item(200, 4).NumberFormat = $ws.Cells.Item(201, 4).NumberFormat
$ws.Cells.Item(200, 5).Value = 4
$ws.Cells.Item(200, 6).Value = 100112043
$ws.Cells.Item(200, 7).Value = "Pepino ensalada"
$ws.Cells.Item(200, 8).Value = "Sin especificar"
$ws.Cells.Item(200, 9).Value = "Primera"
$ws.Cells.Item(200, 10).Value = 800
$ws.Cells.Item(200, 11).Value = 7000
$ws.Cells.Item(200, 12).Value = 8000
$ws.Cells.Item(200, 13).Value = 7500
$ws.Cells.Item(200, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(200, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(200, 16).Value = 107
$ws.Cells.Item(200, 17).Value = 70
$ws.Cells.Item(200, 18).Value = "Hortaliza"
